# Generate Report for Handback
#
# The file "1a4d56c5-df28-4653-9095-a908722463a7.md" has now been handed
# back (its translation round-trip is complete / in sync with en-US) for
# both the zh-cn and de-de locales. Update the status report accordingly:
#   - Overview sheet: flip that row's per-locale status text.
#   - Each locale sheet (zh-cn / de-de): flip the Status cell, stamp a
#     real "Latest Handback DateTime", and record the round-tripped
#     target/handback files (with their hyperlinks), mirroring the
#     existing Handoff File / Source File columns.
# The other tracked file ("a12b8e80-...") stays "Ready for handoff" and
# is left untouched.

$wb = $excel.ActiveWorkbook

# cornflowerblue (FF6495ED), same color used by the other hyperlink cells
# in this workbook; Excel's Font.Color is a BGR-packed long.
# (Hyperlinks.Add already underlines the cell, so only the color needs
# to be corrected to match the existing hyperlink cells' style.)
$hyperlinkColor = 15570276

function Style-AsHyperlink($range) {
    $range.Font.Color = $hyperlinkColor
}

# ---------------------------------------------------------------------
# Overview sheet
# ---------------------------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("B2").Value = "Handed back: in sync with en-US"
$overview.Range("C2").Value = "Handed back: in sync with en-US"

# ---------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------
$zhcn = $wb.Worksheets.Item("zh-cn")

$zhcn.Range("C2").Value = "Handed back: in sync with en-US"

$zhcn.Range("F2").Value = "1a4d56c5-df28-4653-9095-a908722463a7.md"
$zhcn.Hyperlinks.Add(
    $zhcn.Range("F2"),
    "https://github.com/OpenLocalizationTest/oltest/blob/a5ea5b9f4fe4e94a6ba29ed731e04de06076347c/e2e/1a4d56c5-df28-4653-9095-a908722463a7.md",
    "",
    "",
    "1a4d56c5-df28-4653-9095-a908722463a7.md"
) | Out-Null
Style-AsHyperlink $zhcn.Range("F2")

$zhcn.Range("G2").Value = "1a4d56c5-df28-4653-9095-a908722463a7.c4f42a81a10cffb31811f5bde29222eb706e78fb.zh-cn.xlf"
$zhcn.Hyperlinks.Add(
    $zhcn.Range("G2"),
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/0416e60c0ae884232158b822f0bd52324cca8928/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/1a4d56c5-df28-4653-9095-a908722463a7.c4f42a81a10cffb31811f5bde29222eb706e78fb.zh-cn.xlf",
    "",
    "",
    "1a4d56c5-df28-4653-9095-a908722463a7.c4f42a81a10cffb31811f5bde29222eb706e78fb.zh-cn.xlf"
) | Out-Null
Style-AsHyperlink $zhcn.Range("G2")

$zhcn.Range("H2").Value = "2016-03-22 04:36:15"

# ---------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------
$dede = $wb.Worksheets.Item("de-de")

$dede.Range("C2").Value = "Handed back: in sync with en-US"

$dede.Range("F2").Value = "1a4d56c5-df28-4653-9095-a908722463a7.md"
$dede.Hyperlinks.Add(
    $dede.Range("F2"),
    "https://github.com/OpenLocalizationTest/oltest/blob/a5ea5b9f4fe4e94a6ba29ed731e04de06076347c/e2e/1a4d56c5-df28-4653-9095-a908722463a7.md",
    "",
    "",
    "1a4d56c5-df28-4653-9095-a908722463a7.md"
) | Out-Null
Style-AsHyperlink $dede.Range("F2")

$dede.Range("G2").Value = "1a4d56c5-df28-4653-9095-a908722463a7.c4f42a81a10cffb31811f5bde29222eb706e78fb.de-de.xlf"
$dede.Hyperlinks.Add(
    $dede.Range("G2"),
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/bac884cec7d3893b6d263cfeac77ade5a0c4a93d/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/1a4d56c5-df28-4653-9095-a908722463a7.c4f42a81a10cffb31811f5bde29222eb706e78fb.de-de.xlf",
    "",
    "",
    "1a4d56c5-df28-4653-9095-a908722463a7.c4f42a81a10cffb31811f5bde29222eb706e78fb.de-de.xlf"
) | Out-Null
Style-AsHyperlink $dede.Range("G2")

$dede.Range("H2").Value = "2016-03-22 04:36:21"
